# Weekly update: a new daily price record for Orégano (Mercado Mayorista Lo
# Valledor de Santiago) was collected, so it gets inserted as a new row at
# the top of this block (row 366), pushing the existing rows 366-386 down
# to 367-387 and extending the used range from A1:R386 to A1:R387.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 366 (shifts 366..386 -> 367..387)
$ws.Rows.Item(366).Insert()

# Populate the new row 366 with the latest observation
$ws.Range("A366").Value = 6
$ws.Range("B366").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C366").Value = "Metropolitana"
$ws.Range("D366").Value = 45267
$ws.Range("E366").Value = 13
$ws.Range("F366").Value = 100112029
$ws.Range("G366").Value = "Orégano"
$ws.Range("H366").Value = "Sin especificar"
$ws.Range("I366").Value = "Primera"
$ws.Range("J366").Value = 33
$ws.Range("K366").Value = 16000
$ws.Range("L366").Value = 16000
$ws.Range("M366").Value = 16000
$ws.Range("N366").Value = "$/docena de atados"
$ws.Range("O366").Value = "Región Metropolitana"
$ws.Range("P366").Value = 5333
$ws.Range("Q366").Value = 3
$ws.Range("R366").Value = "Hortaliza"
